$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the AutoFilter criteria (TRUE / date filters) while keeping the filter
# dropdowns and the $A$1:$M$211 range, and unhide every data row that the
# previous filter state had hidden.
$ws.ShowAllData()

# Refresh the unique-date count (K) and max-date (L) columns for the rows that
# were already passing the filter (K gets the count, L gets a real date value
# formatted as yyyy-mm-dd instead of a plain text date).
$ws.Range("K3").Value = 13.0
$ws.Range("L3").Value = 45743.0
$ws.Range("L3").NumberFormat = "yyyy-mm-dd"
$ws.Range("K19").Value = 3.0
$ws.Range("L19").Value = 45744.0
$ws.Range("L19").NumberFormat = "yyyy-mm-dd"
$ws.Range("K20").Value = 3.0
$ws.Range("L20").Value = 45743.0
$ws.Range("L20").NumberFormat = "yyyy-mm-dd"
$ws.Range("K44").Value = 8.0
$ws.Range("L44").Value = 45743.0
$ws.Range("L44").NumberFormat = "yyyy-mm-dd"
$ws.Range("K48").Value = 8.0
$ws.Range("L48").Value = 45744.0
$ws.Range("L48").NumberFormat = "yyyy-mm-dd"
$ws.Range("K53").Value = 3.0
$ws.Range("L53").Value = 45743.0
$ws.Range("L53").NumberFormat = "yyyy-mm-dd"
$ws.Range("K64").Value = 3.0
$ws.Range("L64").Value = 45743.0
$ws.Range("L64").NumberFormat = "yyyy-mm-dd"
$ws.Range("K85").Value = 4.0
$ws.Range("L85").Value = 45743.0
$ws.Range("L85").NumberFormat = "yyyy-mm-dd"
$ws.Range("K90").Value = 3.0
$ws.Range("L90").Value = 45744.0
$ws.Range("L90").NumberFormat = "yyyy-mm-dd"
$ws.Range("K118").Value = 3.0
$ws.Range("L118").Value = 45744.0
$ws.Range("L118").NumberFormat = "yyyy-mm-dd"
$ws.Range("K123").Value = 7.0
$ws.Range("L123").Value = 45744.0
$ws.Range("L123").NumberFormat = "yyyy-mm-dd"
$ws.Range("K125").Value = 3.0
$ws.Range("L125").Value = 45744.0
$ws.Range("L125").NumberFormat = "yyyy-mm-dd"
$ws.Range("K132").Value = 4.0
$ws.Range("L132").Value = 45743.0
$ws.Range("L132").NumberFormat = "yyyy-mm-dd"
$ws.Range("K143").Value = 3.0
$ws.Range("L143").Value = 45743.0
$ws.Range("L143").NumberFormat = "yyyy-mm-dd"
$ws.Range("K157").Value = 13.0
$ws.Range("L157").Value = 45743.0
$ws.Range("L157").NumberFormat = "yyyy-mm-dd"
$ws.Range("K165").Value = 3.0
$ws.Range("L165").Value = 45744.0
$ws.Range("L165").NumberFormat = "yyyy-mm-dd"
$ws.Range("K166").Value = 3.0
$ws.Range("L166").Value = 45744.0
$ws.Range("L166").NumberFormat = "yyyy-mm-dd"
$ws.Range("K169").Value = 4.0
$ws.Range("L169").Value = 45744.0
$ws.Range("L169").NumberFormat = "yyyy-mm-dd"
$ws.Range("K170").Value = 7.0
$ws.Range("L170").Value = 45743.0
$ws.Range("L170").NumberFormat = "yyyy-mm-dd"
$ws.Range("K172").Value = 4.0
$ws.Range("L172").Value = 45744.0
$ws.Range("L172").NumberFormat = "yyyy-mm-dd"
$ws.Range("K174").Value = 2.0
$ws.Range("L174").Value = 45743.0
$ws.Range("L174").NumberFormat = "yyyy-mm-dd"
$ws.Range("K176").Value = 3.0
$ws.Range("L176").Value = 45743.0
$ws.Range("L176").NumberFormat = "yyyy-mm-dd"
$ws.Range("K177").Value = 3.0
$ws.Range("L177").Value = 45743.0
$ws.Range("L177").NumberFormat = "yyyy-mm-dd"
$ws.Range("K178").Value = 3.0
$ws.Range("L178").Value = 45744.0
$ws.Range("L178").NumberFormat = "yyyy-mm-dd"
$ws.Range("K187").Value = 3.0
$ws.Range("L187").Value = 45743.0
$ws.Range("L187").NumberFormat = "yyyy-mm-dd"
$ws.Range("K211").Value = 4.0
$ws.Range("L211").Value = 45744.0
$ws.Range("L211").NumberFormat = "yyyy-mm-dd"
